# Update Receptor expression values (and their downstream derived-specificity /
# edge-weight columns) for the Rln3-Rxfp4 sheet, matching the newly recomputed
# TPM-based NATMI output ("update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Receptor avg/total expression changed -> recompute dependents ---
$ws.Range("M2").Value = 5.002662333333333
$ws.Range("N2").Value = 15.007987

$ws.Range("O2").Value = 0.3527593062265809
$ws.Range("P2").Value = 0.3527593062265809
$ws.Range("Q2").Value = 0.571429105025
$ws.Range("R2").Value = 5.142861945225
$ws.Range("S2").Value = 0.3527593062265809
$ws.Range("T2").Value = 0.3527593062265809

# --- Row 3: receptor specificity shifts because the row-2 total changed ---
$ws.Range("O3").Value = 0.3625979570169652
$ws.Range("P3").Value = 0.3625979570169652
$ws.Range("S3").Value = 0.3625979570169652
$ws.Range("T3").Value = 0.3625979570169652

# --- Row 4: same recomputation as row 3 ---
$ws.Range("O4").Value = 0.2846427367564539
$ws.Range("P4").Value = 0.2846427367564539
$ws.Range("S4").Value = 0.2846427367564539
$ws.Range("T4").Value = 0.2846427367564539
